$d = $word.ActiveDocument

# The document currently ends with a trailing empty paragraph (the
# section's last paragraph mark). We need to insert, just before it:
#   1. an empty paragraph
#   2. a new "Heading 1" (Titre1) paragraph, numbered with the same
#      list (numId=1) as the other section headings
#   3. an italic "File name = ..." paragraph
#   4. a plain paragraph of body text

$finalPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$fr = $finalPara.Range
$fr.Collapse(1)
$fr.InsertParagraphBefore()

# --- New Titre1 heading paragraph ---
$finalPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$fr = $finalPara.Range
$fr.Collapse(1)
$fr.InsertParagraphBefore()

$headingPara = $d.Paragraphs.Item($d.Paragraphs.Count - 1)
$headingPara.Style = "Heading 1"
$headingPara.Range.Text = "Recommendations for Integrating a P300-Based Brain" + [char]0x2013 + "Computer Interface in Virtual Reality Environments for Gaming: An Update"

# Reuse the same numbered list (numId=1) as the document's other
# Titre1 headings, continuing the numbering instead of starting a new list.
$firstHeading = $d.Paragraphs.Item(1)
$listTemplate = $firstHeading.Range.ListFormat.ListTemplate
$headingPara.Range.ListFormat.ApplyListTemplate($listTemplate, $true)

# --- New italic "File name = ..." paragraph ---
$finalPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$fr = $finalPara.Range
$fr.Collapse(1)
$fr.InsertParagraphBefore()

$fileNamePara = $d.Paragraphs.Item($d.Paragraphs.Count - 1)
$fileNamePara.Range.Text = "File name = computers-09-00092-v2.pdf"
$fileNamePara.Range.Font.Italic = $true

# --- New plain body paragraph ---
$finalPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$fr = $finalPara.Range
$fr.Collapse(1)
$fr.InsertParagraphBefore()

$bodyPara = $d.Paragraphs.Item($d.Paragraphs.Count - 1)
$bodyPara.Range.Text = "Principal limitation of using EEG for gaming : low transfer rate (preventing movement while using VR)"
